$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-360). Update the value from 45192 (2023-09-23) to
# 45202 (2023-10-03) for all of them, leaving formatting untouched.
$range = $ws.Range("C2:C360")
$range.Value = 45202
